$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Pictures: proofing language ru-RU -> en-US (both w:val and w:eastAsia)
#    Applies to the two inline pictures in the body and the anchored picture
#    in the "first page" header.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shpRange = $d.InlineShapes.Item($i).Range
    $shpRange.LanguageID = "en-US"
    $shpRange.LanguageIDFarEast = "en-US"
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)
    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $hdr = $section.Headers.Item($h)
        if ($hdr.Exists -and $hdr.Shapes.Count -gt 0) {
            $hdrRange = $hdr.Range
            $hdrRange.LanguageID = "en-US"
            $hdrRange.LanguageIDFarEast = "en-US"
        }
    }
}

# ---------------------------------------------------------------------------
# 2. " не принадлежит" -> " не принадлежат", expressed as two runs:
#    " не принадлежа" + "т" (both occurrences in the body).
# ---------------------------------------------------------------------------
function Fix-NePrinadlezit($startSearchPos) {
    $searchRange = $d.Range($startSearchPos, $d.Content.End)
    $found = $searchRange.Duplicate
    $ok = $found.Find.Execute(" не принадлежит", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return -1
    }

    $full = $d.Range($found.Start, $found.End)
    $full.Text = " не принадлежат"
    $newEnd = $full.End

    # Re-establish the preceding run ("K portfelu") as its own run by
    # toggling (and restoring) Bold on it -- this forces a run split
    # without altering the final resolved formatting.
    $precedingRange = $d.Range($full.Start - 10, $full.Start)
    $precedingRange.Bold = 0
    $precedingRange.Bold = 1

    # Split the trailing "t" off into its own run the same way.
    $lastChar = $d.Range($newEnd - 1, $newEnd)
    $lastChar.Bold = 0
    $lastChar.Bold = 1

    return $newEnd
}

$pos = Fix-NePrinadlezit 0
$pos = Fix-NePrinadlezit $pos
